$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cell F1 ("RXNO_DEF") - copy the header style from an
#    existing header cell so it keeps the bold/bordered look.
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "RXNO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# 2. Remove every existing hyperlink (their targets are stale once the
#    rows below are rewritten) - deleting any one range's Hyperlinks
#    collection clears the whole sheet's hyperlink set in this host.
# ---------------------------------------------------------------------------
$ws.Range("B2").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 3. Rewrite the data rows (2-9) with the updated mapping content.
#    Columns: A index, B BAO_IRI, C BAO_DESC, D RXNO_IRI, E RXNO_DESC,
#    F RXNO_DEF
# ---------------------------------------------------------------------------
$rows = @(
  @{ Row=2; A=0; B="http://purl.obolibrary.org/obo/CHEBI_23367"; C="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_23367'}"; D="http://purl.obolibrary.org/obo/CHEBI_23367"; E="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_23367'}"; F="['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']" },
  @{ Row=3; A=1; B="http://purl.obolibrary.org/obo/CHEBI_39141"; C="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"; D="http://purl.obolibrary.org/obo/CHEBI_39141"; E="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39141'}"; F="[]" },
  @{ Row=4; A=2; B="http://purl.obolibrary.org/obo/CHEBI_39142"; C="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39142'}"; D="http://purl.obolibrary.org/obo/CHEBI_39142"; E="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39142'}"; F="[]" },
  @{ Row=5; A=3; B="http://purl.obolibrary.org/obo/CHEBI_39143"; C="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39143'}"; D="http://purl.obolibrary.org/obo/CHEBI_39143"; E="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39143'}"; F="[]" },
  @{ Row=6; A=4; B="http://purl.obolibrary.org/obo/CHEBI_39144"; C="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39144'}"; D="http://purl.obolibrary.org/obo/CHEBI_39144"; E="{'iri': 'http://purl.obolibrary.org/obo/CHEBI_39144'}"; F="[]" },
  @{ Row=7; A=5; B="http://purl.obolibrary.org/obo/CHEBI_52214"; C="{'label': 'ligand', 'prefLabel': None, 'altLabel': None, 'name': 'CHEBI_52214'}"; D="http://purl.obolibrary.org/obo/MOP_0000714"; E="{'label': 'ligand'}"; F="[]" },
  @{ Row=8; A=6; B="http://www.bioassayontology.org/bao#BAO_0003043"; C="{'label': 'molecular entity', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0003043'}"; D="http://purl.obolibrary.org/obo/CHEBI_23367"; E="{'label': 'molecular entity', 'prefLabel': 'molecular entity'}"; F="['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']" },
  @{ Row=9; A=7; B="http://www.bioassayontology.org/bao#BAO_0000325"; C="{'label': 'peptide', 'prefLabel': None, 'altLabel': None, 'name': 'BAO_0000325'}"; D="http://purl.obolibrary.org/obo/CHEBI_16670"; E="{'label': 'peptide'}"; F="[]" }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("A$n").Value = $r.A
  # Rows 8/9 are brand-new rows - give column A the same bold/bordered
  # "index" look the existing rows (2-7) already carry.
  if ($n -ge 8) {
    $ws.Range("A$n").Font.Bold = $true
    $ws.Range("A$n").Borders.LineStyle = 1
    $ws.Range("A$n").HorizontalAlignment = -4108
    $ws.Range("A$n").VerticalAlignment = -4160
  }
  $ws.Range("B$n").Value = $r.B
  $ws.Range("B$n").Style = "Hyperlink"
  $ws.Range("C$n").Value = $r.C
  $ws.Range("D$n").Value = $r.D
  $ws.Range("D$n").Style = "Hyperlink"
  $ws.Range("E$n").Value = $r.E
  $ws.Range("F$n").Value = $r.F
}

# ---------------------------------------------------------------------------
# 4. Re-create the hyperlinks in order (B then D for each row, 2 through 9)
#    matching the refreshed cell contents. Rows 8/9's B hyperlink keeps the
#    BAO ontology fragment anchor.
# ---------------------------------------------------------------------------
$links = @(
  @{ Cell="B2"; Address="http://purl.obolibrary.org/obo/CHEBI_23367"; SubAddress="" },
  @{ Cell="D2"; Address="http://purl.obolibrary.org/obo/CHEBI_23367"; SubAddress="" },
  @{ Cell="B3"; Address="http://purl.obolibrary.org/obo/CHEBI_39141"; SubAddress="" },
  @{ Cell="D3"; Address="http://purl.obolibrary.org/obo/CHEBI_39141"; SubAddress="" },
  @{ Cell="B4"; Address="http://purl.obolibrary.org/obo/CHEBI_39142"; SubAddress="" },
  @{ Cell="D4"; Address="http://purl.obolibrary.org/obo/CHEBI_39142"; SubAddress="" },
  @{ Cell="B5"; Address="http://purl.obolibrary.org/obo/CHEBI_39143"; SubAddress="" },
  @{ Cell="D5"; Address="http://purl.obolibrary.org/obo/CHEBI_39143"; SubAddress="" },
  @{ Cell="B6"; Address="http://purl.obolibrary.org/obo/CHEBI_39144"; SubAddress="" },
  @{ Cell="D6"; Address="http://purl.obolibrary.org/obo/CHEBI_39144"; SubAddress="" },
  @{ Cell="B7"; Address="http://purl.obolibrary.org/obo/CHEBI_52214"; SubAddress="" },
  @{ Cell="D7"; Address="http://purl.obolibrary.org/obo/MOP_0000714"; SubAddress="" },
  @{ Cell="B8"; Address="http://www.bioassayontology.org/bao"; SubAddress="BAO_0003043" },
  @{ Cell="D8"; Address="http://purl.obolibrary.org/obo/CHEBI_23367"; SubAddress="" },
  @{ Cell="B9"; Address="http://www.bioassayontology.org/bao"; SubAddress="BAO_0000325" },
  @{ Cell="D9"; Address="http://purl.obolibrary.org/obo/CHEBI_16670"; SubAddress="" }
)

foreach ($lk in $links) {
  if ($lk.SubAddress -ne "") {
    $ws.Hyperlinks.Add($ws.Range($lk.Cell), $lk.Address, $lk.SubAddress)
  } else {
    $ws.Hyperlinks.Add($ws.Range($lk.Cell), $lk.Address)
  }
}

Write-Host "Done. Hyperlinks: $($ws.Hyperlinks.Count)"
